$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.671.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.349.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.62%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.658"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.57"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +14.81%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +19.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0983"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "27.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.47%  "
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "17.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +15.14%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.54%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.699.16"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +10.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.880"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.361.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.630.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("E19").Value = "  +4.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "75.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "250.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.24%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.81%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.64%  "
$ws.Range("E27").Value = "  -2.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "172.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.74%  "
$ws.Range("E30").Value = "  +8.18%  "
$ws.Range("E31").Value = "  +2.79%  "
$ws.Range("E32").Value = "  +4.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0700"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.12"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.76"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.72%  "
$ws.Range("E37").Value = "  +7.92%  "
$ws.Range("E38").Value = "  +1.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0264"
$ws.Range("D39").Style = "Normal"
$ws.Range("E40").Value = "  +13.63%  "
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("E43").Value = "  +8.98%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.76%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "98.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0964"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("E48").Value = "  +13.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.439.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.10%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.19%  "
$ws.Range("B51").Value = "TerraClassic"
$ws.Range("C51").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000204"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.17%  "
